# This script applies the "gh-pages output generated at 456a3b4" update to
# 杭州-漫展信息.xlsx: a handful of "想去人数" (interest-count) bumps across
# the 展览 / 演出 / 本地生活 / 全部类型 sheets, one ticket that sold out
# (numeric price -> "已售罄"), and the 全部类型 (combined) sheet catching up
# with two events that were already present in 展览/演出 but missing from
# the combined view (杭州·乙游Camerata杭州2.0 and 杭州·黄西全新脱口秀专场
# 《水土不服》), which pushes a couple of neighbouring rows' content down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 1483
$ws1.Range("F5").Value = 3235
$ws1.Range("F7").Value = 653
$ws1.Range("F8").Value = 2213
$ws1.Range("F9").Value = 476
$ws1.Range("F14").Value = 1065
$ws1.Range("F16").Value = 5
$ws1.Range("F18").Value = 199
$ws1.Range("F19").Value = 4418
$ws1.Range("F20").Value = 1287
$ws1.Range("F21").Value = 3364
$ws1.Range("F23").Value = 66
$ws1.Range("F24").Value = 159
$ws1.Range("F25").Value = 3290
$ws1.Range("F26").Value = 4896
$ws1.Range("F30").Value = 3179
$ws1.Range("F31").Value = 342
$ws1.Range("F36").Value = 1146
$ws1.Range("G36").Value = "已售罄"
$ws1.Range("F37").Value = 1390
$ws1.Range("F38").Value = 112
$ws1.Range("F39").Value = 1318
$ws1.Range("F40").Value = 840
$ws1.Range("F41").Value = 12
$ws1.Range("F45").Value = 282
$ws1.Range("F46").Value = 58
$ws1.Range("F48").Value = 362
$ws1.Range("F49").Value = 3708

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F20").Value = 48
$ws2.Range("F23").Value = 14

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value = 2090

# ---------------------------------------------------------------------
# Sheet "全部类型" (All categories - combined view)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 1483
$ws4.Range("F5").Value = 3235
$ws4.Range("F7").Value = 653
$ws4.Range("F9").Value = 2213
$ws4.Range("F10").Value = 476
$ws4.Range("F16").Value = 1065

# Row 18 becomes the (previously missing) 乙游Camerata event.
# (leading apostrophe forces text so Excel doesn't auto-convert the
# yyyy-mm-dd-looking string into a real date, matching the other date
# cells in this column which are plain text too)
$ws4.Range("B18").Value = "'2024-07-04"
$ws4.Range("C18").Value = "杭州·乙游Camerata杭州2.0"
$ws4.Range("D18").Value = "杭海路601号江和美海洋广场1层 嘉宝丽酒店"
$ws4.Range("E18").Value = "2024.07.04 10:00-07.04 17:00"
$ws4.Range("F18").Value = 5
$ws4.Range("G18").Value = 68
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=86777"
$ws4.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202406/aBDjuHlA1717403033570.jpeg"

# Row 19 becomes what used to be in row 18 (黑执事only), with refreshed counts.
$ws4.Range("B19").Value = "'2024-07-06"
$ws4.Range("C19").Value = "杭州·黑执事only"
$ws4.Range("D19").Value = "大岭山路156号 爱丽芬城堡"
$ws4.Range("E19").Value = "2024.07.06 10:00-07.07 18:00"
$ws4.Range("F19").Value = 199
$ws4.Range("G19").Value = 160
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=86414"
$ws4.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202405/iP2cxk2w1716800288950.jpeg"

# Row 20 becomes what used to be in row 19 (AD04动漫展), with refreshed counts.
$ws4.Range("C20").Value = "杭州·AD04动漫展"
$ws4.Range("D20").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws4.Range("E20").Value = "2024.07.13 10:00-07.14 17:00"
$ws4.Range("F20").Value = 4418
$ws4.Range("G20").Value = 75
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=85012"
$ws4.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202405/y1iKqqnh1715326769523.jpeg"

$ws4.Range("F21").Value = 1287
$ws4.Range("F23").Value = 3364
$ws4.Range("F24").Value = 3290
$ws4.Range("F25").Value = 4896
$ws4.Range("F28").Value = 3179
$ws4.Range("F29").Value = 342

# Row 34 becomes the (previously missing) 黄西 stand-up show.
$ws4.Range("C34").Value = "杭州·黄西全新脱口秀专场《水土不服》"
$ws4.Range("D34").Value = "延安路279号 浙江胜利剧院"
$ws4.Range("E34").Value = "2024.07.27 19:30-07.27 21:30"
$ws4.Range("F34").Value = 2
$ws4.Range("G34").Value = 224
$ws4.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=85935"
$ws4.Range("I34").Value = "//i2.hdslb.com/bfs/openplatform/202405/9YqhR4Ke1716191781870.jpeg"

$ws4.Range("F35").Value = 1390
$ws4.Range("F36").Value = 112
$ws4.Range("F37").Value = 1318
$ws4.Range("F39").Value = 840
$ws4.Range("F41").Value = 48
$ws4.Range("F44").Value = 282
$ws4.Range("F45").Value = 14
$ws4.Range("F46").Value = 58
$ws4.Range("F48").Value = 362
$ws4.Range("F49").Value = 3708
